$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.167.29"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.07"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.43"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5235"
$ws.Range("E7").Value = "  +1.81%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3770"
$ws.Range("E8").Value = "  +0.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07252"
$ws.Range("E9").Value = "  +0.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.21"
$ws.Range("E10").Value = "  -0.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9042"
$ws.Range("E11").Value = "  -0.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08501"
$ws.Range("E12").Value = "  +11.21%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.913.07"
$ws.Range("E13").Value = "  +1.19%  "

# Row 14
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.90"
$ws.Range("E14").Value = "  +2.17%  "

# Row 15
$ws.Range("E15").Value = "  +0.53%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008673"
$ws.Range("E17").Value = "  +2.29%  "

# Row 18
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.202.21"
$ws.Range("E20").Value = "  +0.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.086"
$ws.Range("E21").Value = "  +0.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.144.97"
$ws.Range("E22").Value = "  +1.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  +0.59%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.439"
$ws.Range("E24").Value = "  +0.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.325"
$ws.Range("E25").Value = "  +1.91%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.05"
$ws.Range("E26").Value = "  +0.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.24"
$ws.Range("E27").Value = "  +1.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.745"
$ws.Range("E28").Value = "  -1.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.13"
$ws.Range("E29").Value = "  +0.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.931"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.821"
$ws.Range("E31").Value = "  -0.18%  "

# Row 32
$ws.Range("E32").Value = "  +1.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8056"
$ws.Range("E33").Value = "  +3.08%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05057"
$ws.Range("E34").Value = "  -0.67%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.247"
$ws.Range("E35").Value = "  +0.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.451"
$ws.Range("E36").Value = "  +5.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.951"
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.608"
$ws.Range("E38").Value = "  -0.56%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5735"
$ws.Range("E39").Value = "  +2.62%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02000"
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.138"
$ws.Range("E42").Value = "  +0.36%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.637"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "115.98"
$ws.Range("E44").Value = "  -1.56%  "

# Row 45
$ws.Range("E45").Value = "  +0.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4869"
$ws.Range("E46").Value = "  +1.24%  "

# Row 47
$ws.Range("E47").Value = "  -0.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.0000"
$ws.Range("E48").Value = "  +0.10%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.622"
$ws.Range("E49").Value = "  +1.30%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.72"
$ws.Range("E50").Value = "  +0.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.23"
$ws.Range("E51").Value = "  +0.35%  "
